# Update the crypto price/volume table (rows 2-51) with refreshed figures.
# Numeric-looking price strings are entered with a leading apostrophe so
# Excel keeps them as literal text (matching the original formatting,
# e.g. "227.81" rather than being reinterpreted as the number 227.81).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.879.35"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "2.034.28"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'227.81"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").Value = "'60.35"
$ws.Range("E7").Value = "  +3.81%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.387"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("E10").Value = "  +1.06%  "

$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").Value = "'14.69"
$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").Value = "2.332.79"
$ws.Range("E13").Value = "  -1.96%  "

$ws.Range("D14").Value = "'21.19"
$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").Value = "'0.757"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").Value = "'5.24"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").Value = "2.048.08"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").Value = "37.826.01"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").Value = "'6.04"
$ws.Range("E19").Value = "  -3.42%  "

$ws.Range("D20").Value = "'69.87"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("E21").Value = "  -0.93%  "

$ws.Range("D22").Value = "'225.91"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("E24").Value = "  -2.06%  "

$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").Value = "'9.27"
$ws.Range("E26").Value = "  -0.42%  "

$ws.Range("D27").Value = "'164.93"
$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("E28").Value = "  -3.93%  "

$ws.Range("D29").Value = "'18.91"
$ws.Range("E29").Value = "  -1.15%  "

$ws.Range("E30").Value = "  -6.73%  "

$ws.Range("E31").Value = "  +1.41%  "

$ws.Range("D32").Value = "'4.44"
$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("E33").Value = "  +3.49%  "

$ws.Range("D34").Value = "'0.0603"
$ws.Range("E34").Value = "  -2.21%  "

$ws.Range("D35").Value = "'4.48"
$ws.Range("E35").Value = "  -2.63%  "

$ws.Range("D36").Value = "'6.40"
$ws.Range("E36").Value = "  +5.28%  "

$ws.Range("D37").Value = "'2.26"
$ws.Range("E37").Value = "  -5.76%  "

$ws.Range("D38").Value = "'3.26"
$ws.Range("E38").Value = "  -1.16%  "

$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").Value = "1.537.95"
$ws.Range("E40").Value = "  +3.56%  "

$ws.Range("D41").Value = "'0.0218"
$ws.Range("E41").Value = "  -0.79%  "

$ws.Range("D42").Value = "'96.97"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").Value = "'16.75"
$ws.Range("E43").Value = "  -0.47%  "

$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("D45").Value = "'0.0922"
$ws.Range("E45").Value = "  -3.03%  "

$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("D47").Value = "'3.96"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  -2.20%  "

$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("D50").Value = "'7.14"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("D51").Value = "2.222.56"
$ws.Range("E51").Value = "  -1.85%  "
